$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.461.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.192.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.61"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.23"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.60"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.78"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.518.22"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.51"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.200.08"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.410.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.12"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.02"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.42"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.13"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.28"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.121"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0769"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.24"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.83"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.56"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0306"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.76"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.37%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.55"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0992"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.36"
